# Add a new "Estimates" worksheet after the existing "Sheet1" and populate
# it with a story-point estimate summary row.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so the tab order is Sheet1, Estimates.
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Estimates"

# Column widths / row height to match the source layout (values chosen so the
# COM "characters" width rounds to the same stored width as the original file).
$ws.Columns.Item(2).ColumnWidth = 55.1666666666667
$ws.Columns.Item(3).ColumnWidth = 9.70963541666667
$ws.Rows.Item(2).RowHeight = 43.5

# Content.
$ws.Range("B2").Value = "Total Story Points Esimates (Including Desing, Cut Effort, DB Design, Testing, Requirement Detailing, Code Review, Bug Fixing, Documentation, Release Notes)"
$ws.Range("B2").WrapText = $true
$ws.Range("C2").Value = 314

# Leave the new sheet as the active sheet/tab with C2 selected.
$ws.Range("C2").Select() | Out-Null
